$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12208
$ws1.Range("F3").Value = 6868
$ws1.Range("F5").Value = 20
$ws1.Range("F6").Value = 434
$ws1.Range("F8").Value = 11
$ws1.Range("F9").Value = 13
$ws1.Range("F10").Value = 956
$ws1.Range("F11").Value = 113
$ws1.Range("F12").Value = 312
$ws1.Range("F13").Value = 962
$ws1.Range("F14").Value = 3692
$ws1.Range("F16").Value = 997
$ws1.Range("F17").Value = 506
$ws1.Range("F18").Value = 211
$ws1.Range("F19").Value = 329
$ws1.Range("F21").Value = 246
$ws1.Range("F22").Value = 275
$ws1.Range("F25").Value = 5106
$ws1.Range("F27").Value = 1295
$ws1.Range("F28").Value = 261
$ws1.Range("F29").Value = 770
$ws1.Range("F30").Value = 1260

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 88
$ws2.Range("F4").Value = 3718
$ws2.Range("F13").Value = 7

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9165
$ws3.Range("F3").Value = 532
$ws3.Range("F4").Value = 1903

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9165
$ws4.Range("F3").Value = 532
$ws4.Range("F4").Value = 1903
$ws4.Range("F5").Value = 12208
$ws4.Range("F6").Value = 6868
$ws4.Range("F7").Value = 88
$ws4.Range("F8").Value = 3718
$ws4.Range("F10").Value = 20
$ws4.Range("F11").Value = 434
$ws4.Range("F13").Value = 11
$ws4.Range("F14").Value = 13
$ws4.Range("F15").Value = 956
$ws4.Range("F16").Value = 113
$ws4.Range("F17").Value = 312
$ws4.Range("F18").Value = 962
$ws4.Range("F19").Value = 3692
$ws4.Range("F21").Value = 997
$ws4.Range("F22").Value = 506
$ws4.Range("F23").Value = 211
$ws4.Range("F24").Value = 329
$ws4.Range("F26").Value = 246
$ws4.Range("F27").Value = 275
$ws4.Range("F33").Value = 5106
$ws4.Range("F35").Value = 1295
$ws4.Range("F38").Value = 261
$ws4.Range("F40").Value = 770
$ws4.Range("F41").Value = 1260
$ws4.Range("F43").Value = 7
